$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (Row 9) ---
$ws.Range("S9").Value = "ready to be hadded"
$ws.Range("T9").Value = "ready to be hadded"

# --- Cell value updates (Row 10) ---
$ws.Range("K10").Value = "ready to be fit"
$ws.Range("S10").Value = "ready to be hadded"
$ws.Range("T10").Value = "ready to be hadded"

# --- Cell value updates (Row 11) ---
$ws.Range("C11").Value = "ready to be fit"
$ws.Range("H11").Value = "ready to be fit"
$ws.Range("J11").Value = "looks like 12 jobs didn't run properly"
$ws.Range("K11").Value = "ready to be hadded"
$ws.Range("S11").Value = "ready to be hadded"
$ws.Range("T11").Value = "looks like many jobs are missing"

# --- Column width updates ---
# ColumnWidth is in "characters"; Excel stores a slightly different
# effective width (characters + ~0.8333) in the xlsx <col> element, so we
# compensate here to land exactly on the target stored widths (22 / 28).
$ws.Columns.Item(19).ColumnWidth = 21.16666666666667
$ws.Columns.Item(20).ColumnWidth = 27.16666666666667

# --- Selection update ---
$ws.Range("J15").Select()
